# Lagt til alternativ for ingen vakanse.
# Updates "Etterspørsel" (col C) and "Vakanse" (col E) figures for the
# five education-category blocks on Sheet1, reflecting the new
# "no vacancy" alternative scenario used in the teacher-demand model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=3; C=46514; E=1491},
    @{Row=4; C=45964; E=3090},
    @{Row=5; C=45811; E=4300},
    @{Row=6; C=45543; E=5630},
    @{Row=7; C=45395; E=6836},
    @{Row=8; C=45335; E=7955},
    @{Row=9; C=45373; E=8949},
    @{Row=10; C=45497; E=9830},
    @{Row=11; C=45701; E=10572},
    @{Row=12; C=45938; E=11204},
    @{Row=13; C=46188; E=11747},
    @{Row=14; C=46436; E=12192},
    @{Row=15; C=46690; E=12536},
    @{Row=16; C=46981; E=12766},
    @{Row=17; C=47303; E=12878},
    @{Row=18; C=47647; E=12921},
    @{Row=19; C=48009; E=12926},
    @{Row=20; C=48343; E=12968},
    @{Row=21; C=48604; E=13113},
    @{Row=22; C=48789; E=13385},
    @{Row=24; C=52344; E=428},
    @{Row=25; C=52374; E=857},
    @{Row=26; C=52235; E=1511},
    @{Row=27; C=51975; E=2342},
    @{Row=28; C=51531; E=3395},
    @{Row=29; C=50990; E=4551},
    @{Row=30; C=50499; E=5687},
    @{Row=31; C=50066; E=6744},
    @{Row=32; C=49630; E=7767},
    @{Row=33; C=49310; E=8609},
    @{Row=34; C=49033; E=9349},
    @{Row=35; C=48813; E=9954},
    @{Row=36; C=48652; E=10420},
    @{Row=37; C=48651; E=10669},
    @{Row=38; C=48756; E=10754},
    @{Row=39; C=48927; E=10751},
    @{Row=40; C=49144; E=10693},
    @{Row=41; C=49409; E=10583},
    @{Row=42; C=49728; E=10442},
    @{Row=43; C=50065; E=10298},
    @{Row=46; C=16643; E=2003},
    @{Row=47; C=16688; E=3056},
    @{Row=48; C=16715; E=4130},
    @{Row=49; C=16747; E=5199},
    @{Row=50; C=16769; E=6259},
    @{Row=51; C=16769; E=7329},
    @{Row=52; C=16717; E=8423},
    @{Row=53; C=16644; E=9523},
    @{Row=54; C=16569; E=10592},
    @{Row=55; C=16516; E=11604},
    @{Row=56; C=16478; E=12576},
    @{Row=57; C=16465; E=13502},
    @{Row=58; C=16427; E=14450},
    @{Row=59; C=16377; E=15395},
    @{Row=60; C=16340; E=16308},
    @{Row=61; C=16347; E=17172},
    @{Row=62; C=16380; E=18013},
    @{Row=63; C=16419; E=18844},
    @{Row=64; C=16468; E=19649},
    @{Row=66; C=30773; E=20},
    @{Row=67; C=30884; E=6},
    @{Row=68; C=30964; E=40},
    @{Row=69; C=31009; E=120},
    @{Row=70; C=31068; E=188},
    @{Row=71; C=31106; E=277},
    @{Row=72; C=31101; E=394},
    @{Row=73; C=30987; E=606},
    @{Row=74; C=30830; E=841},
    @{Row=75; C=30668; E=1049},
    @{Row=76; C=30548; E=1177},
    @{Row=77; C=30457; E=1241},
    @{Row=78; C=30415; E=1217},
    @{Row=79; C=30317; E=1235},
    @{Row=80; C=30190; E=1246},
    @{Row=81; C=30088; E=1233},
    @{Row=82; C=30078; E=1112},
    @{Row=83; C=30118; E=950},
    @{Row=84; C=30171; E=775},
    @{Row=85; C=30244; E=590},
    @{Row=87; C=13277; E=-52},
    @{Row=88; C=13341; E=-150},
    @{Row=89; C=13402; E=-262},
    @{Row=90; C=13451; E=-360},
    @{Row=91; C=13522; E=-498},
    @{Row=92; C=13589; E=-644},
    @{Row=93; C=13624; E=-759},
    @{Row=94; C=13588; E=-825},
    @{Row=95; C=13527; E=-863},
    @{Row=96; C=13453; E=-897},
    @{Row=97; C=13401; E=-967},
    @{Row=98; C=13363; E=-1043},
    @{Row=99; C=13351; E=-1160},
    @{Row=100; C=13294; E=-1214},
    @{Row=101; C=13210; E=-1249},
    @{Row=102; C=13137; E=-1285},
    @{Row=103; C=13115; E=-1367},
    @{Row=104; C=13120; E=-1486},
    @{Row=105; C=13129; E=-1588},
    @{Row=106; C=13149; E=-1706}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
